$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples")

# Reorder columns: move specId (A) to become column D, shifting
# builder/model/generator left. This also carries the column widths
# along, matching a real drag-and-drop column move in Excel.
$ws.Columns.Item(1).Cut()
$ws.Columns.Item(5).Insert()

# New rows describing additional webapp samples.
$ws.Range("A4").Value2 = "PredefinedWebAppBuilder"
$ws.Range("B4").Value2 = "WebAppModel"
$ws.Range("C4").Value2 = "WebAppGenerator"
$ws.Range("D4").Value2 = "webappAngular"
$ws.Range("E4").Value2 = "Wildfly, UserResource, tutti, war"

$ws.Range("A5").Value2 = "PredefinedWebAppBuilder"
$ws.Range("B5").Value2 = "WebAppModel"
$ws.Range("C5").Value2 = "WebAppGenerator"
$ws.Range("D5").Value2 = "webappFirst"
$ws.Range("E5").Value2 = "Tomcat, UserResource"

# New columns F/G on the spring4RestTomcat row.
$ws.Range("F3").Value2 = "buildAppSimpleSpring"
$ws.Range("G3").Value2 = "+"
$ws.Columns.Item(6).ColumnWidth = 18.5
$ws.Columns.Item(7).ColumnWidth = 2.1666666666666665

# Sort rows 2:5 by columns A,B,C,D ascending (matches the workbook's
# recorded sortState).
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A5"))
$ws.Sort.SortFields.Add($ws.Range("B2:B5"))
$ws.Sort.SortFields.Add($ws.Range("C2:C5"))
$ws.Sort.SortFields.Add($ws.Range("D2:D5"))
$ws.Sort.SetRange($ws.Range("A2:E5"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

$ws.Range("F4").Select() | Out-Null

$ws1 = $wb.Worksheets.Item("classes")
$ws1.Range("A1").Select() | Out-Null
$ws.Activate() | Out-Null
